# Turns the single empty "Sheet1" workbook into a small German<->English
# vocabulary workbook with three sheets:
#   1) "Hallo und guten tag"
#   2) "Beruf und Familie"
#   3) "Beruf und Familie sen"
#
# Cell-entry order below mirrors the order the shared strings appear in the
# target file (sheet2's A1:A3/B1:B2 first, then all of sheet3, then back to
# sheet2 for B3/K10) so the shared-strings table indices line up.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet 1: "Hallo und guten tag" -----------------------------------
$ws1.Range("A2").Value = "Hallo"
$ws1.Range("B2").Value = "hello"

# --- Sheet 2: "Beruf und Familie" (inserted after sheet 1) ------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Range("A1").Value = "Frage "
$ws2.Range("B1").Value = "question"
$ws2.Range("A2").Value = "Antwort "
$ws2.Range("B2").Value = "answer "
$ws2.Range("A3").Value = "Wohnort"

# --- Sheet 3: "Beruf und Familie sen" (inserted after sheet 2) --------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Range("A1").Value = "We heißen Sie? "
$ws3.Range("B1").Value = "What is your name? "
$ws3.Range("A2").Value = "Ich heiße... "
$ws3.Range("A3").Value = "Mein Name ist…"
$ws3.Range("B2").Value = "My name is... "
$ws3.Range("B3").Value = "My name is... "

# Finish off sheet 2's remaining cells.
$ws2.Range("B3").Value = "place of residence"
$ws2.Range("K10").Value = " "

# --- Rename sheets ------------------------------------------------------
$ws1.Name = "Hallo und guten tag"
$ws2.Name = "Beruf und Familie"
$ws3.Name = "Beruf und Familie sen"

# --- Recreate the saved selections / active sheet -----------------------
$ws1.Range("G27").Select()
$ws3.Range("A4").Select()
$ws2.Activate()
$ws2.Range("K10").Select()
